$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) is unchanged. Rows 2-7 are rebuilt as a 3x2
# sending/target-cluster combination matrix (ECs, FAPs, sCs) for the
# Ngf -> Ntrk1 interaction, replacing the single previous data row.

$data = @(
    @("ECs",  "Ngf", "Ntrk1", "FAPs", 1, 0.3333333333333333, 0.2143436666666667, 0.643031, 0.02290233455071953, 0.02290233455071953, 3, 1, 0.1211523333333333, 0.363457, 0.6524744903454663, 0.6524744903454663, 0.02596823535188889, 0.233714118167, 0.01494318906370209, 0.01494318906370209),
    @("ECs",  "Ngf", "Ntrk1", "sCs",  1, 0.3333333333333333, 0.2143436666666667, 0.643031, 0.02290233455071953, 0.02290233455071953, 1, 0.3333333333333333, 0.064529, 0.193587, 0.3475255096545335, 0.3475255096545335, 0.01383138246633333, 0.124482442197, 0.007959145487017438, 0.007959145487017438),
    @("FAPs", "Ngf", "Ntrk1", "FAPs", 2, 0.6666666666666666, 0.666149, 1.998447, 0.07117713108058832, 0.07117713108058833, 3, 1, 0.1211523333333333, 0.363457, 0.6524744903454663, 0.6524744903454663, 0.08070550569766667, 0.7263495512790001, 0.04644126232605932, 0.04644126232605932),
    @("FAPs", "Ngf", "Ntrk1", "sCs",  2, 0.6666666666666666, 0.666149, 1.998447, 0.07117713108058832, 0.07117713108058833, 1, 0.3333333333333333, 0.064529, 0.193587, 0.3475255096545335, 0.3475255096545335, 0.042985928821, 0.386873359389, 0.024735868754529, 0.024735868754529),
    @("sCs",  "Ngf", "Ntrk1", "FAPs", 3, 1, 8.478538666666667, 25.435616, 0.9059205343686921, 0.9059205343686921, 3, 1, 0.1211523333333333, 0.363457, 0.6524744903454663, 0.6524744903454663, 1.027194742723556, 9.244752684512001, 0.5910900389557049, 0.5910900389557049),
    @("sCs",  "Ngf", "Ntrk1", "sCs",  3, 1, 8.478538666666667, 25.435616, 0.9059205343686921, 0.9059205343686921, 1, 0.3333333333333333, 0.064529, 0.193587, 0.3475255096545335, 0.3475255096545335, 0.5471116216213334, 4.924004594592001, 0.3148304954129871, 0.3148304954129871)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
